$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '39.524.82'
$ws.Range("E2").Value = '  +1.85%  '
$ws.Range("D3").Value = '2.169.61'
$ws.Range("E3").Value = '  +3.55%  '
$ws.Range("D5").Value = '''230.33'
$ws.Range("E5").Value = '  +0.90%  '
$ws.Range("D6").Value = '''0.624'
$ws.Range("E6").Value = '  +1.46%  '
$ws.Range("D7").Value = '''65.03'
$ws.Range("E7").Value = '  +7.03%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("D9").Value = '''0.400'
$ws.Range("E9").Value = '  +3.71%  '
$ws.Range("D10").Value = '''0.0870'
$ws.Range("E10").Value = '  +2.87%  '
$ws.Range("D11").Value = '''0.104'
$ws.Range("E11").Value = '  +0.08%  '
$ws.Range("D12").Value = '''16.08'
$ws.Range("E12").Value = '  +6.17%  '
$ws.Range("D13").Value = '2.491.66'
$ws.Range("E13").Value = '  +3.59%  '
$ws.Range("D14").Value = '''22.51'
$ws.Range("E14").Value = '  +2.57%  '
$ws.Range("D15").Value = '''0.822'
$ws.Range("E15").Value = '  +1.62%  '
$ws.Range("D16").Value = '''5.61'
$ws.Range("E16").Value = '  +2.33%  '
$ws.Range("D17").Value = '2.165.09'
$ws.Range("E17").Value = '  +3.59%  '
$ws.Range("D18").Value = '39.522.28'
$ws.Range("E18").Value = '  +1.90%  '
$ws.Range("D19").Value = '''72.69'
$ws.Range("E19").Value = '  +1.33%  '
$ws.Range("E20").Value = '  +0.71%  '
$ws.Range("D21").Value = '0.0₃0858'
$ws.Range("E21").Value = '  +1.87%  '
$ws.Range("D22").Value = '''232.45'
$ws.Range("E22").Value = '  +2.21%  '
$ws.Range("E23").Value = '  -0.02%  '
$ws.Range("B24").Value = 'Toncoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D24").Value = '''2.37'
$ws.Range("E24").Value = '  -0.27%  '
$ws.Range("B25").Value = 'PancakeSwap'
$ws.Range("C25").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D25").Value = '''2.39'
$ws.Range("E25").Value = '  +2.41%  '
$ws.Range("D26").Value = '''9.74'
$ws.Range("E26").Value = '  +1.80%  '
$ws.Range("D27").Value = '''172.60'
$ws.Range("E27").Value = '  +0.71%  '
$ws.Range("D28").Value = '''0.139'
$ws.Range("E28").Value = '  -1.10%  '
$ws.Range("D29").Value = '''20.09'
$ws.Range("E29").Value = '  +4.40%  '
$ws.Range("D30").Value = '''1.42'
$ws.Range("E30").Value = '  -2.52%  '
$ws.Range("D31").Value = '''2.74'
$ws.Range("E31").Value = '  +12.54%  '
$ws.Range("D32").Value = '''0.123'
$ws.Range("E32").Value = '  +2.25%  '
$ws.Range("D33").Value = '''4.69'
$ws.Range("E33").Value = '  +3.35%  '
$ws.Range("E34").Value = '  +2.85%  '
$ws.Range("D35").Value = '''7.15'
$ws.Range("E35").Value = '  +9.73%  '
$ws.Range("D36").Value = '''0.0625'
$ws.Range("E36").Value = '  +2.52%  '
$ws.Range("D37").Value = '''2.44'
$ws.Range("E37").Value = '  +1.62%  '
$ws.Range("D38").Value = '''3.61'
$ws.Range("E38").Value = '  +0.15%  '
$ws.Range("D39").Value = '''0.999'
$ws.Range("E39").Value = '  -0.13%  '
$ws.Range("D40").Value = '''105.13'
$ws.Range("E40").Value = '  +4.37%  '
$ws.Range("D41").Value = '''0.0232'
$ws.Range("E41").Value = '  +0.74%  '
$ws.Range("D42").Value = '''17.99'
$ws.Range("E42").Value = '  +0.00%  '
$ws.Range("D43").Value = '1.537.32'
$ws.Range("E43").Value = '  -0.22%  '
$ws.Range("E44").Value = '  +6.55%  '
$ws.Range("D45").Value = '''0.0932'
$ws.Range("E45").Value = '  +1.47%  '
$ws.Range("E46").Value = '  +7.28%  '
$ws.Range("E47").Value = '  -0.01%  '
$ws.Range("D48").Value = '''7.87'
$ws.Range("E48").Value = '  +2.00%  '
$ws.Range("D49").Value = '''4.26'
$ws.Range("E49").Value = '  +2.99%  '
$ws.Range("D50").Value = '2.374.36'
$ws.Range("E50").Value = '  +3.52%  '
$ws.Range("D51").Value = '''2.98'
$ws.Range("E51").Value = '  +0.24%  '
